$d = $word.ActiveDocument

# ------------------------------------------------------------------
# Change 1: {{ p / lots / HeaderLabel }}  ->  {{ plotsHeaderLabel }}
# (three runs merge into a single run, keeping the first run's formatting)
# ------------------------------------------------------------------
$null = $d.Content.Find.Execute(
    "{{ p" + "lots" + "HeaderLabel }}",
    $false, $false, $false, $false, $false, $true, 1, $false,
    "{{ plotsHeaderLabel }}", 2)

# ------------------------------------------------------------------
# Change 2: {% if p / lot / sHeader%}{{ p / lots / Header }}{% else %}-{% endif%}
#   -> {% if plotsHeader%}{{ plotsHeader }}{% else %}-{% endif%}
# (five runs merge into a single run, keeping the first run's formatting)
# ------------------------------------------------------------------
$null = $d.Content.Find.Execute(
    "{% if p" + "lot" + "sHeader%}{{ p" + "lots" + "Header }}{% else %}-{% endif%}",
    $false, $false, $false, $false, $false, $true, 1, $false,
    "{% if plotsHeader%}{{ plotsHeader }}{% else %}-{% endif%}", 2)

# ------------------------------------------------------------------
# Change 3: the inputDateHeader table cell.
#   - w:ilvl goes from 5 to 2
#   - the single run's text is split into three runs so a new
#     paperInputDateHeader fallback is introduced between the
#     original {{ inputDateHeader }} value and the else/endif tail.
# ------------------------------------------------------------------
$r = $d.Content
$null = $r.Find.Execute(
    "{% if inputDateHeader %}{{ inputDateHeader }}{% else %}-{% endif%}",
    $false, $false, $false, $false, $false, $true, 1, $false, "", 0)

$xml = '<?xml version="1.0" standalone="yes"?>' +
'<?mso-application progid="Word.Document"?>' +
'<pkg:package xmlns:pkg="http://schemas.microsoft.com/office/2006/xmlPackage">' +
'<pkg:part pkg:name="/word/document.xml" pkg:contentType="application/vnd.openxmlformats-officedocument.wordprocessingml.document.main+xml">' +
'<pkg:xmlData>' +
'<w:document xmlns:w="http://schemas.openxmlformats.org/wordprocessingml/2006/main">' +
'<w:body>' +
'<w:p>' +
  '<w:pPr>' +
    '<w:pStyle w:val="Heading3"/>' +
    '<w:numPr><w:ilvl w:val="2"/><w:numId w:val="3"/></w:numPr>' +
    '<w:suppressAutoHyphens w:val="true"/>' +
    '<w:spacing w:before="0" w:after="0"/>' +
    '<w:ind w:right="176" w:hanging="0"/>' +
    '<w:rPr>' +
      '<w:rFonts w:ascii="Liberation Sans" w:hAnsi="Liberation Sans" w:eastAsia="Noto Sans CJK SC Regular" w:cs="Arial"/>' +
      '<w:b w:val="false"/><w:bCs w:val="false"/>' +
      '<w:color w:val="auto"/>' +
      '<w:kern w:val="2"/>' +
      '<w:sz w:val="18"/><w:szCs w:val="18"/>' +
      '<w:lang w:val="fr-CH" w:eastAsia="zh-CN" w:bidi="hi-IN"/>' +
    '</w:rPr>' +
  '</w:pPr>' +
  '<w:r>' +
    '<w:rPr>' +
      '<w:rFonts w:eastAsia="Noto Sans CJK SC Regular" w:cs="Arial" w:ascii="Liberation Sans" w:hAnsi="Liberation Sans"/>' +
      '<w:b w:val="false"/><w:bCs w:val="false"/>' +
      '<w:color w:val="auto"/>' +
      '<w:kern w:val="2"/>' +
      '<w:sz w:val="18"/><w:szCs w:val="18"/>' +
      '<w:lang w:val="fr-CH" w:eastAsia="zh-CN" w:bidi="hi-IN"/>' +
    '</w:rPr>' +
    '<w:t>{% if inputDateHeader %}{{ inputDateHeader }}</w:t>' +
  '</w:r>' +
  '<w:r>' +
    '<w:rPr>' +
      '<w:rFonts w:eastAsia="Noto Sans CJK SC Regular" w:cs="Arial" w:ascii="Liberation Sans" w:hAnsi="Liberation Sans"/>' +
      '<w:b w:val="false"/><w:bCs w:val="false"/>' +
      '<w:color w:val="auto"/>' +
      '<w:kern w:val="2"/>' +
      '<w:sz w:val="18"/><w:szCs w:val="18"/>' +
      '<w:lang w:val="fr-CH" w:eastAsia="zh-CN" w:bidi="hi-IN"/>' +
    '</w:rPr>' +
    '<w:t>{% if paperInputDateHeader %}({{ paperInputDateHeader }}){% else %}{% endif %}</w:t>' +
  '</w:r>' +
  '<w:r>' +
    '<w:rPr>' +
      '<w:rFonts w:eastAsia="Noto Sans CJK SC Regular" w:cs="Arial" w:ascii="Liberation Sans" w:hAnsi="Liberation Sans"/>' +
      '<w:b w:val="false"/><w:bCs w:val="false"/>' +
      '<w:color w:val="auto"/>' +
      '<w:kern w:val="2"/>' +
      '<w:sz w:val="18"/><w:szCs w:val="18"/>' +
      '<w:lang w:val="fr-CH" w:eastAsia="zh-CN" w:bidi="hi-IN"/>' +
    '</w:rPr>' +
    '<w:t>{% else %}-{% endif %}</w:t>' +
  '</w:r>' +
'</w:p>' +
'</w:body>' +
'</w:document>' +
'</pkg:xmlData></pkg:part></pkg:package>'

$null = $r.InsertXML($xml)
